# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
# These reflect a re-run of the scrape job at a later wall-clock time;
# every other column is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-16T07:01:17.075648+00:00"
$ws.Range("K3").Value = "2025-12-16T07:01:17.075682+00:00"
$ws.Range("K4").Value = "2025-12-16T07:01:17.075702+00:00"
$ws.Range("K5").Value = "2025-12-16T07:01:19.366260+00:00"
$ws.Range("K6").Value = "2025-12-16T07:01:19.366295+00:00"
$ws.Range("K7").Value = "2025-12-16T07:01:19.366313+00:00"
$ws.Range("K8").Value = "2025-12-16T07:01:22.144126+00:00"
$ws.Range("K9").Value = "2025-12-16T07:01:26.179669+00:00"
$ws.Range("K10").Value = "2025-12-16T07:01:29.026042+00:00"
$ws.Range("K11").Value = "2025-12-16T07:01:31.751573+00:00"
$ws.Range("K12").Value = "2025-12-16T07:01:36.382237+00:00"
$ws.Range("K13").Value = "2025-12-16T07:01:36.382264+00:00"
$ws.Range("K14").Value = "2025-12-16T07:01:38.672017+00:00"
$ws.Range("K15").Value = "2025-12-16T07:01:41.613560+00:00"
$ws.Range("K16").Value = "2025-12-16T07:01:44.422151+00:00"
$ws.Range("K17").Value = "2025-12-16T07:01:47.165636+00:00"
$ws.Range("K18").Value = "2025-12-16T07:01:47.165670+00:00"
$ws.Range("K19").Value = "2025-12-16T07:01:47.165689+00:00"
$ws.Range("K20").Value = "2025-12-16T07:01:47.165705+00:00"
$ws.Range("K21").Value = "2025-12-16T07:01:47.165721+00:00"
$ws.Range("K22").Value = "2025-12-16T07:01:49.494849+00:00"
$ws.Range("K23").Value = "2025-12-16T07:01:49.494880+00:00"
$ws.Range("K24").Value = "2025-12-16T07:01:52.414950+00:00"
$ws.Range("K25").Value = "2025-12-16T07:01:52.414978+00:00"
$ws.Range("K26").Value = "2025-12-16T07:01:52.414995+00:00"
$ws.Range("K27").Value = "2025-12-16T07:01:52.415011+00:00"
$ws.Range("K28").Value = "2025-12-16T07:01:52.415025+00:00"
$ws.Range("K29").Value = "2025-12-16T07:01:55.231924+00:00"
$ws.Range("K30").Value = "2025-12-16T07:01:55.231953+00:00"
$ws.Range("K31").Value = "2025-12-16T07:01:55.231971+00:00"
$ws.Range("K32").Value = "2025-12-16T07:01:58.099316+00:00"
$ws.Range("K33").Value = "2025-12-16T07:01:58.099346+00:00"
$ws.Range("K34").Value = "2025-12-16T07:01:58.099363+00:00"
$ws.Range("K35").Value = "2025-12-16T07:02:00.535533+00:00"
$ws.Range("K36").Value = "2025-12-16T07:02:02.986738+00:00"
$ws.Range("K37").Value = "2025-12-16T07:02:02.986767+00:00"
$ws.Range("K38").Value = "2025-12-16T07:02:08.627813+00:00"
$ws.Range("K39").Value = "2025-12-16T07:02:08.627840+00:00"
$ws.Range("K40").Value = "2025-12-16T07:02:11.360455+00:00"
$ws.Range("K41").Value = "2025-12-16T07:02:11.360482+00:00"
